$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A (Indice_Punteo) so all other columns shift left by one.
$ws.Range("A:A").Delete()

# Add two new columns at the end (now Q and R) both titled "Indice_Punteo",
# copying the header style from the other header cells.
$ws.Range("Q1").Value = "Indice_Punteo"
$ws.Range("R1").Value = "Indice_Punteo"
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "Indice_Punteo"
$ws.Range("R1").Value = "Indice_Punteo"

# New trailing cells in row 2 stay blank for the two new columns, but still
# materialize as real (empty) cells, matching the source workbook.
$ws.Range("Q2").Style = "Normal"
$ws.Range("R2").Style = "Normal"
